$wb = $excel.ActiveWorkbook

# --- 1. Add the new "metadata" worksheet, positioned right after "data" ---
$wsData = $wb.Worksheets.Item("data")
$wsNew = $wb.Worksheets.Add()
$wsNew.Name = "metadata"
$wsNew.Move($null, $wb.Worksheets.Item("data"))

# Re-resolve references by name (Move()/Copy() invalidate earlier handles).
$wsData = $wb.Worksheets.Item("data")
$ws2 = $wb.Worksheets.Item("metadata")

# --- 2. Copy the header styling (bold/border/center) from "data" onto the
#        new header row + the leading index-column cell ---
$wsData.Range("B1:F1").Copy()
$ws2.Range("B1:F1").PasteSpecial(-4122)
$wsData = $wb.Worksheets.Item("data")
$ws2 = $wb.Worksheets.Item("metadata")

$wsData.Range("F1").Copy()
$ws2.Range("G1").PasteSpecial(-4122)
$wsData = $wb.Worksheets.Item("data")
$ws2 = $wb.Worksheets.Item("metadata")

$wsData.Range("A2").Copy()
$ws2.Range("A2").PasteSpecial(-4122)
$wsData = $wb.Worksheets.Item("data")
$ws2 = $wb.Worksheets.Item("metadata")

# --- 3. Populate the metadata header row ---
$ws2.Range("B1").Value = "data_name"
$ws2.Range("C1").Value = "data_id"
$ws2.Range("D1").Value = "data_version"
$ws2.Range("E1").Value = "data_version_created"
$ws2.Range("F1").Value = "panel_query_time"
$ws2.Range("G1").Value = "panel_get_request"

# --- 4. Populate the single metadata data row ---
$ws2.Range("A2").Value = 0
$ws2.Range("B2").Value = "Familial diabetes"
$ws2.Range("C2").Value = 152

# data_version ("1.62") must land as literal text, not a parsed number -
# format the cell as Text first, write the value, then restore the
# (stock) "Normal" style so no stray "s" attribute is left on the cell.
$ws2.Range("D2").NumberFormat = "@"
$ws2.Range("D2").Value = "1.62"
$ws2.Range("D2").Style = "Normal"
$ws2 = $wb.Worksheets.Item("metadata")

$ws2.Range("E2").Value = "2021-07-28T09:56:25.863987Z"
$ws2.Range("F2").Value = "2021-10-05 14:20:12.418688"
$ws2.Range("G2").Value = "https://panelapp.genomicsengland.co.uk/api/v1/panels/152/?format=json"

# --- 5. Refresh the F-column "time_taken" timestamps on the "data" sheet ---
$wsData = $wb.Worksheets.Item("data")
$times = @(
    "2021-10-05 14:20:12.422347",
    "2021-10-05 14:20:12.422354",
    "2021-10-05 14:20:12.422357",
    "2021-10-05 14:20:12.422360",
    "2021-10-05 14:20:12.422363",
    "2021-10-05 14:20:12.422365",
    "2021-10-05 14:20:12.422368",
    "2021-10-05 14:20:12.422370",
    "2021-10-05 14:20:12.422373",
    "2021-10-05 14:20:12.422375",
    "2021-10-05 14:20:12.422378",
    "2021-10-05 14:20:12.422380",
    "2021-10-05 14:20:12.422383",
    "2021-10-05 14:20:12.422385",
    "2021-10-05 14:20:12.422388",
    "2021-10-05 14:20:12.422390",
    "2021-10-05 14:20:12.422393",
    "2021-10-05 14:20:12.422396",
    "2021-10-05 14:20:12.422398",
    "2021-10-05 14:20:12.422401",
    "2021-10-05 14:20:12.422403",
    "2021-10-05 14:20:12.422405",
    "2021-10-05 14:20:12.422408",
    "2021-10-05 14:20:12.422410",
    "2021-10-05 14:20:12.422413",
    "2021-10-05 14:20:12.422416",
    "2021-10-05 14:20:12.422418",
    "2021-10-05 14:20:12.422420",
    "2021-10-05 14:20:12.422423",
    "2021-10-05 14:20:12.422425",
    "2021-10-05 14:20:12.422428",
    "2021-10-05 14:20:12.422430",
    "2021-10-05 14:20:12.422433",
    "2021-10-05 14:20:12.422435",
    "2021-10-05 14:20:12.422438",
    "2021-10-05 14:20:12.422440",
    "2021-10-05 14:20:12.422443",
    "2021-10-05 14:20:12.422445",
    "2021-10-05 14:20:12.422448",
    "2021-10-05 14:20:12.422450",
    "2021-10-05 14:20:12.422453",
    "2021-10-05 14:20:12.422455",
    "2021-10-05 14:20:12.422458",
    "2021-10-05 14:20:12.422460",
    "2021-10-05 14:20:12.422463",
    "2021-10-05 14:20:12.422465",
    "2021-10-05 14:20:12.422468",
    "2021-10-05 14:20:12.422470",
    "2021-10-05 14:20:12.422472",
    "2021-10-05 14:20:12.422475",
    "2021-10-05 14:20:12.422477",
    "2021-10-05 14:20:12.422480",
    "2021-10-05 14:20:12.422483",
    "2021-10-05 14:20:12.422485",
    "2021-10-05 14:20:12.422488",
    "2021-10-05 14:20:12.422490"
)
for ($i = 0; $i -lt $times.Length; $i++) {
    $row = $i + 2
    $wsData.Range("F$row").Value = $times[$i]
}

Write-Output "edit complete"
